# Auto-generated cell updates for the "Chocobo_Profits" workbook scheduled-runner sync.
# Applies per-row H:N recomputed market-price/profit figures across the 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 6000
$ws.Range("M31").Value = -5770
$ws.Range("H40").Value = 4120.4287
$ws.Range("I40").Value = 2326.6667
$ws.Range("J40").Value = 4609.636
$ws.Range("K40").Value = 2326.6667
$ws.Range("L40").Value = 4609.636
$ws.Range("M40").Value = -2151.6667
$ws.Range("N40").Value = -4959.636
$ws.Range("H41").Value = 350.05884
$ws.Range("I41").Value = 270.63635
$ws.Range("J41").Value = 495.66666
$ws.Range("K41").Value = 270.63635
$ws.Range("L41").Value = 495.66666
$ws.Range("M41").Value = 169.36365
$ws.Range("N41").Value = -1375.66666
$ws.Range("H112").Value = 1346.8292
$ws.Range("J112").Value = 1360.5
$ws.Range("L112").Value = 4081.5
$ws.Range("N112").Value = -6297.5
$ws.Range("H113").Value = 6545.231
$ws.Range("I113").Value = 1747.5
$ws.Range("J113").Value = 7417.5454
$ws.Range("K113").Value = 1747.5
$ws.Range("L113").Value = 7417.5454
$ws.Range("M113").Value = 1506.5
$ws.Range("N113").Value = -13925.5454
$ws.Range("H132").Value = 33035072
$ws.Range("I132").Value = 39002410
$ws.Range("J132").Value = 2004933.2
$ws.Range("K132").Value = 117007230
$ws.Range("L132").Value = 6014799.6
$ws.Range("M132").Value = -117004700
$ws.Range("N132").Value = -6019859.6
$ws.Range("H135").Value = 994.7273
$ws.Range("I135").Value = 990.75
$ws.Range("J135").Value = 1005.3333
$ws.Range("K135").Value = 8916.75
$ws.Range("L135").Value = 9047.9997
$ws.Range("M135").Value = -6381.75
$ws.Range("N135").Value = -14117.9997
$ws.Range("H137").Value = 810453.1
$ws.Range("I137").Value = 2074257
$ws.Range("K137").Value = 6222771
$ws.Range("M137").Value = -6220221
$ws.Range("H138").Value = 2973.2322
$ws.Range("I138").Value = 1647.9445
$ws.Range("J138").Value = 3601
$ws.Range("K138").Value = 4943.833500000001
$ws.Range("L138").Value = 10803
$ws.Range("M138").Value = 196.1664999999994
$ws.Range("N138").Value = -21083
$ws.Range("H141").Value = 19997.666
$ws.Range("I141").Value = 28045.25
$ws.Range("J141").Value = 3902.5
$ws.Range("K141").Value = 84135.75
$ws.Range("L141").Value = 11707.5
$ws.Range("M141").Value = -78955.75
$ws.Range("N141").Value = -22067.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2707.889
$ws.Range("I61").Value = 2910.1428
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2910.1428
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -2698.1428
$ws.Range("N61").Value = -2424
$ws.Range("H109").Value = 30501.291
$ws.Range("J109").Value = 30501.291
$ws.Range("L109").Value = 30501.291
$ws.Range("N109").Value = -33275.291
$ws.Range("H110").Value = 640.4583
$ws.Range("I110").Value = 506.375
$ws.Range("J110").Value = 908.625
$ws.Range("K110").Value = 506.375
$ws.Range("L110").Value = 908.625
$ws.Range("M110").Value = 1538.625
$ws.Range("N110").Value = -4998.625
$ws.Range("H132").Value = 2350.0278
$ws.Range("I132").Value = 1961.2069
$ws.Range("J132").Value = 3960.8572
$ws.Range("K132").Value = 5883.620699999999
$ws.Range("L132").Value = 11882.5716
$ws.Range("M132").Value = -3353.620699999999
$ws.Range("N132").Value = -16942.5716
$ws.Range("H136").Value = 2707.889
$ws.Range("I136").Value = 2910.1428
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 8730.428400000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -6180.428400000001
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 57500
$ws.Range("J58").Value = 57500
$ws.Range("L58").Value = 57500
$ws.Range("N58").Value = -58088
$ws.Range("H133").Value = 43774.5
$ws.Range("J133").Value = 56840
$ws.Range("L133").Value = 56840
$ws.Range("N133").Value = -66960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H31").Value = 235411.34
$ws.Range("I31").Value = 795918.3
$ws.Range("J31").Value = 3006.0244
$ws.Range("K31").Value = 795918.3
$ws.Range("L31").Value = 3006.0244
$ws.Range("M31").Value = -795623.3
$ws.Range("N31").Value = -3596.0244
$ws.Range("H34").Value = 235411.34
$ws.Range("I34").Value = 795918.3
$ws.Range("J34").Value = 3006.0244
$ws.Range("K34").Value = 795918.3
$ws.Range("L34").Value = 3006.0244
$ws.Range("M34").Value = -795716.3
$ws.Range("N34").Value = -3410.0244
$ws.Range("H132").Value = 8005.3
$ws.Range("I132").Value = 9400
$ws.Range("J132").Value = 6610.6
$ws.Range("K132").Value = 28200
$ws.Range("L132").Value = 19831.8
$ws.Range("M132").Value = -25670
$ws.Range("N132").Value = -24891.8
$ws.Range("H134").Value = 9142.857
$ws.Range("I134").Value = 10078.454
$ws.Range("J134").Value = 5712.3335
$ws.Range("K134").Value = 30235.362
$ws.Range("L134").Value = 17137.0005
$ws.Range("M134").Value = -27700.362
$ws.Range("N134").Value = -22207.0005
$ws.Range("H137").Value = 44542.855
$ws.Range("J137").Value = 44542.855
$ws.Range("L137").Value = 44542.855
$ws.Range("N137").Value = -54742.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2529.4285
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H68").Value = 2958.1333
$ws.Range("J68").Value = 3706.6511
$ws.Range("L68").Value = 11119.9533
$ws.Range("N68").Value = -12741.9533
$ws.Range("H71").Value = 2958.1333
$ws.Range("J71").Value = 3706.6511
$ws.Range("L71").Value = 33359.8599
$ws.Range("N71").Value = -41471.8599
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H131").Value = 839.6767599999999
$ws.Range("I131").Value = 515.5
$ws.Range("J131").Value = 860.5914
$ws.Range("K131").Value = 1546.5
$ws.Range("L131").Value = 2581.7742
$ws.Range("M131").Value = 3493.5
$ws.Range("N131").Value = -12661.7742
$ws.Range("H132").Value = 2811.3845
$ws.Range("I132").Value = 933.3333
$ws.Range("J132").Value = 3374.8
$ws.Range("K132").Value = 8399.9997
$ws.Range("L132").Value = 30373.2
$ws.Range("M132").Value = -5869.9997
$ws.Range("N132").Value = -35433.2
$ws.Range("H137").Value = 2505.4614
$ws.Range("I137").Value = 2413.9
$ws.Range("J137").Value = 2810.6667
$ws.Range("K137").Value = 7241.700000000001
$ws.Range("L137").Value = 8432.000100000001
$ws.Range("M137").Value = -2141.700000000001
$ws.Range("N137").Value = -18632.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6495.7437
$ws.Range("I70").Value = 5843.3105
$ws.Range("J70").Value = 8387.799999999999
$ws.Range("K70").Value = 5843.3105
$ws.Range("L70").Value = 8387.799999999999
$ws.Range("M70").Value = -5573.3105
$ws.Range("N70").Value = -8927.799999999999
$ws.Range("H73").Value = 6495.7437
$ws.Range("I73").Value = 5843.3105
$ws.Range("J73").Value = 8387.799999999999
$ws.Range("K73").Value = 5843.3105
$ws.Range("L73").Value = 8387.799999999999
$ws.Range("M73").Value = -4907.3105
$ws.Range("N73").Value = -10259.8
$ws.Range("H80").Value = 2847.3635
$ws.Range("I80").Value = 2703.125
$ws.Range("K80").Value = 2703.125
$ws.Range("M80").Value = -1705.125
$ws.Range("H83").Value = 2847.3635
$ws.Range("I83").Value = 2703.125
$ws.Range("K83").Value = 13515.625
$ws.Range("M83").Value = -8523.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 41100.32
$ws.Range("I22").Value = 72069.92999999999
$ws.Range("J22").Value = 1684.4546
$ws.Range("K22").Value = 72069.92999999999
$ws.Range("L22").Value = 1684.4546
$ws.Range("M22").Value = -71774.92999999999
$ws.Range("N22").Value = -2274.4546
$ws.Range("H27").Value = 41100.32
$ws.Range("I27").Value = 72069.92999999999
$ws.Range("J27").Value = 1684.4546
$ws.Range("K27").Value = 72069.92999999999
$ws.Range("L27").Value = 1684.4546
$ws.Range("M27").Value = -71962.92999999999
$ws.Range("N27").Value = -1898.4546
$ws.Range("H68").Value = 1069.8689
$ws.Range("I68").Value = 987.8305
$ws.Range("K68").Value = 987.8305
$ws.Range("M68").Value = -238.8305
$ws.Range("H71").Value = 1069.8689
$ws.Range("I71").Value = 987.8305
$ws.Range("K71").Value = 4939.1525
$ws.Range("M71").Value = -1195.1525
$ws.Range("H122").Value = 7333.3335
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 9750
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 29250
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -34150
$ws.Range("H132").Value = 5154.095
$ws.Range("I132").Value = 4041.25
$ws.Range("J132").Value = 5838.923
$ws.Range("K132").Value = 12123.75
$ws.Range("L132").Value = 17516.769
$ws.Range("M132").Value = -9593.75
$ws.Range("N132").Value = -22576.769
$ws.Range("H136").Value = 2956.5454
$ws.Range("I136").Value = 1247.5238
$ws.Range("K136").Value = 3742.5714
$ws.Range("M136").Value = -1192.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2883.1538
$ws.Range("I132").Value = 1275.8889
$ws.Range("K132").Value = 3827.6667
$ws.Range("M132").Value = -1297.6667
$ws.Range("H136").Value = 2603.7173
$ws.Range("I136").Value = 1090.2333
$ws.Range("J136").Value = 5441.5
$ws.Range("K136").Value = 3270.699900000001
$ws.Range("L136").Value = 16324.5
$ws.Range("M136").Value = -720.6999000000005
$ws.Range("N136").Value = -21424.5
